# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    7   = @("b",  "Acknowledge (Backchannel)")
    8   = @("aa", "Agree/Accept")
    15  = @("aa", "Agree/Accept")
    18  = @("aa", "Agree/Accept")
    20  = @("sd", "Statement-non-opinion")
    21  = @("sd", "Statement-non-opinion")
    23  = @("sd", "Statement-non-opinion")
    26  = @("sd", "Statement-non-opinion")
    46  = @("sd", "Statement-non-opinion")
    57  = @("sd", "Statement-non-opinion")
    65  = @("b",  "Acknowledge (Backchannel)")
    66  = @("sv", "Statement-opinion")
    72  = @("aa", "Agree/Accept")
    74  = @("sd", "Statement-non-opinion")
    78  = @("ba", "Appreciation")
    93  = @("sd", "Statement-non-opinion")
    124 = @("aa", "Agree/Accept")
    125 = @("sd", "Statement-non-opinion")
    127 = @("sv", "Statement-opinion")
    134 = @("b",  "Acknowledge (Backchannel)")
    143 = @("sd", "Statement-non-opinion")
    178 = @("aa", "Agree/Accept")
    197 = @("%",  "Uninterpretable")
    205 = @("b",  "Acknowledge (Backchannel)")
    207 = @("aa", "Agree/Accept")
    211 = @("sv", "Statement-opinion")
}

foreach ($row in $changes.Keys) {
    $values = $changes[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}
